$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.382629871368408
$ws.Range("B1").Value = 4.421439170837402
$ws.Range("C1").Value = 5.674187183380127
$ws.Range("D1").Value = 8.814920425415039
$ws.Range("E1").Value = 5.765583038330078
